$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($rng, $val)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '42.406.77'
$ws.Range('E2').Value = '  -1.34%  '

Set-TextValue $ws.Range('D3') '2.280.29'
$ws.Range('E3').Value = '  -0.74%  '

Set-TextValue $ws.Range('D5') '304.84'
$ws.Range('E5').Value = '  +1.94%  '

Set-TextValue $ws.Range('D6') '95.34'
$ws.Range('E6').Value = '  -2.26%  '

Set-TextValue $ws.Range('D7') '0.502'
$ws.Range('E7').Value = '  -2.92%  '

$ws.Range('E8').Value = '  -0.01%  '

Set-TextValue $ws.Range('D9') '0.492'
$ws.Range('E9').Value = '  -3.44%  '

Set-TextValue $ws.Range('D10') '34.76'
$ws.Range('E10').Value = '  -3.92%  '

Set-TextValue $ws.Range('D11') '0.0782'
$ws.Range('E11').Value = '  -0.65%  '

$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range('D12') '0.118'
$ws.Range('E12').Value = '  +1.07%  '

$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D13') '17.93'
$ws.Range('E13').Value = '  +0.28%  '

Set-TextValue $ws.Range('D14') '6.64'
$ws.Range('E14').Value = '  -2.07%  '

Set-TextValue $ws.Range('D15') '2.639.32'
$ws.Range('E15').Value = '  -0.52%  '

Set-TextValue $ws.Range('D16') '2.285.95'
$ws.Range('E16').Value = '  -1.15%  '

Set-TextValue $ws.Range('D17') '0.770'
$ws.Range('E17').Value = '  -1.70%  '

Set-TextValue $ws.Range('D18') '42.343.65'
$ws.Range('E18').Value = '  -1.34%  '

Set-TextValue $ws.Range('D19') '12.60'
$ws.Range('E19').Value = '  -1.11%  '

Set-TextValue $ws.Range('D20') '0.0₃0886'
$ws.Range('E20').Value = '  -2.36%  '

Set-TextValue $ws.Range('D21') '5.98'
$ws.Range('E21').Value = '  -1.96%  '

Set-TextValue $ws.Range('D22') '66.89'
$ws.Range('E22').Value = '  -2.81%  '

Set-TextValue $ws.Range('D23') '234.79'
$ws.Range('E23').Value = '  -2.95%  '

Set-TextValue $ws.Range('D24') '2.15'
$ws.Range('E24').Value = '  -0.74%  '

$ws.Range('E25').Value = '  +0.00%  '

Set-TextValue $ws.Range('D26') '2.42'
$ws.Range('E26').Value = '  -0.05%  '

Set-TextValue $ws.Range('D27') '24.73'
$ws.Range('E27').Value = '  -0.20%  '

Set-TextValue $ws.Range('D28') '165.61'
$ws.Range('E28').Value = '  +0.04%  '

$ws.Range('E29').Value = '  +0.51%  '

Set-TextValue $ws.Range('D30') '8.93'
$ws.Range('E30').Value = '  -1.40%  '

Set-TextValue $ws.Range('D31') '32.25'
$ws.Range('E31').Value = '  -2.20%  '

$ws.Range('E32').Value = '  +0.08%  '

Set-TextValue $ws.Range('D33') '4.91'
$ws.Range('E33').Value = '  -2.09%  '

$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range('D34') '17.42'
$ws.Range('E34').Value = '  -1.60%  '

$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D35') '4.58'
$ws.Range('E35').Value = '  -3.24%  '

Set-TextValue $ws.Range('D36') '2.36'
$ws.Range('E36').Value = '  -1.85%  '

Set-TextValue $ws.Range('D37') '0.0681'
$ws.Range('E37').Value = '  -1.15%  '

Set-TextValue $ws.Range('D38') '0.100'
$ws.Range('E38').Value = '  -1.45%  '

Set-TextValue $ws.Range('D39') '1.72'
$ws.Range('E39').Value = '  -1.94%  '

$ws.Range('E40').Value = '  -2.30%  '

Set-TextValue $ws.Range('D41') '2.66'
$ws.Range('E41').Value = '  -3.37%  '

Set-TextValue $ws.Range('D42') '1.990.02'
$ws.Range('E42').Value = '  -0.97%  '

Set-TextValue $ws.Range('D43') '0.0275'
$ws.Range('E43').Value = '  -3.20%  '

$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D44') '9.94'
$ws.Range('E44').Value = '  -2.36%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D45') '17.72'
$ws.Range('E45').Value = '  +3.01%  '

Set-TextValue $ws.Range('D46') '2.02'
$ws.Range('E46').Value = '  -7.60%  '

Set-TextValue $ws.Range('D47') '2.74'
$ws.Range('E47').Value = '  -1.98%  '

$ws.Range('E48').Value = '  +8.45%  '

Set-TextValue $ws.Range('D49') '53.08'
$ws.Range('E49').Value = '  -1.52%  '

Set-TextValue $ws.Range('D50') '2.506.19'
$ws.Range('E50').Value = '  -0.47%  '

Set-TextValue $ws.Range('D51') '70.68'
$ws.Range('E51').Value = '  -3.07%  '
